# The workbook has a single sheet with a table (Tabella1) listing funicular
# stations. The first table/header column was renamed from "station" to
# "name" (the rest of the header row - address/town/province - stays the
# same). Renaming the header cell updates both the worksheet cell and the
# table's column name (ListObject), since the header cell IS the table
# column header.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "name"
